$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update fluid_mass (B3) which feeds the Q_toBoil formula in B7
$ws.Range("B3").Value = 0.58874952899999999

$wb.Save()
